$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price strings; some (e.g. "6.50") look like plain
# numbers, and a direct .Value assignment would make Excel
# auto-convert them to floating point Number cells (e.g. 6.5099999999999998)
# and change the cells type/format. To keep them as text (matching the
# original inline-string cells) without touching any cell style, we push
# the new text through a temporary ="..." formula and then paste its
# result back in as a literal value.
function Set-TextValue($addr, $value) {
    $escaped = $value -replace '"', '""'
    $ws.Range($addr).Formula = '="' + $escaped + '"'
    $ws.Range($addr).Copy()
    $ws.Range($addr).PasteSpecial(-4163)
    $excel.CutCopyMode = $false
}

Set-TextValue "D2" '68.063.03'
$ws.Range("E2").Value = '  -0.91%  '

Set-TextValue "D3" '3.775.45'
$ws.Range("E3").Value = '  -2.39%  '

$ws.Range("E4").Value = '  +0.04%  '

Set-TextValue "D5" '596.37'
$ws.Range("E5").Value = '  -1.03%  '

Set-TextValue "D6" '168.56'
$ws.Range("E6").Value = '  -1.90%  '

Set-TextValue "D7" '3.775.88'

$ws.Range("E8").Value = '  +0.02%  '

$ws.Range("E9").Value = '  -0.78%  '

$ws.Range("E10").Value = '  -3.15%  '

Set-TextValue "D11" '6.51'
$ws.Range("E11").Value = '  +0.16%  '

$ws.Range("E12").Value = '  -2.43%  '

Set-TextValue "D13" '0.0000279'
$ws.Range("E13").Value = '  -3.24%  '

Set-TextValue "D14" '36.61'
$ws.Range("E14").Value = '  -1.67%  '

Set-TextValue "D15" '4.409.04'
$ws.Range("E15").Value = '  -2.23%  '

Set-TextValue "D16" '3.779.24'
$ws.Range("E16").Value = '  -2.01%  '

Set-TextValue "D17" '18.68'
$ws.Range("E17").Value = '  +1.57%  '

Set-TextValue "D18" '68.065.57'
$ws.Range("E18").Value = '  -0.92%  '

Set-TextValue "D19" '7.15'
$ws.Range("E19").Value = '  -3.40%  '

$ws.Range("E20").Value = '  -0.45%  '

$ws.Range("E21").Value = '  -4.90%  '

Set-TextValue "D22" '466.81'
$ws.Range("E22").Value = '  -1.30%  '

Set-TextValue "D23" '0.716'
$ws.Range("E23").Value = '  -2.31%  '

Set-TextValue "D24" '0.0000149'
$ws.Range("E24").Value = '  -8.82%  '

Set-TextValue "D25" '83.85'
$ws.Range("E25").Value = '  -0.07%  '

$ws.Range("E26").Value = '  -1.66%  '

Set-TextValue "D27" '12.13'
$ws.Range("E27").Value = '  -0.95%  '

Set-TextValue "D28" '10.32'
$ws.Range("E28").Value = '  -1.75%  '

$ws.Range("B30").Value = 'PancakeSwap'
$ws.Range("C30").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
Set-TextValue "D30" '2.91'
$ws.Range("E30").Value = '  -1.30%  '

$ws.Range("B31").Value = 'WrappedeETH'
$ws.Range("C31").Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
Set-TextValue "D31" '3.924.04'
$ws.Range("E31").Value = '  -2.33%  '

$ws.Range("E32").Value = '  -2.76%  '

Set-TextValue "D33" '30.40'
$ws.Range("E33").Value = '  -3.16%  '

$ws.Range("E34").Value = '  -4.13%  '

Set-TextValue "D35" '9.20'
$ws.Range("E35").Value = '  -2.35%  '

Set-TextValue "D36" '3.732.89'
$ws.Range("E36").Value = '  -2.52%  '

Set-TextValue "D37" '3.79'
$ws.Range("E37").Value = '  -4.36%  '

$ws.Range("E38").Value = '  -1.91%  '

$ws.Range("E39").Value = '  -1.23%  '

$ws.Range("E40").Value = '  -1.77%  '

$ws.Range("E41").Value = '  -2.56%  '

Set-TextValue "D42" '0.999'
$ws.Range("E42").Value = '  -0.03%  '

$ws.Range("E43").Value = '  -2.99%  '

$ws.Range("E44").Value = '  +0.00%  '

Set-TextValue "D45" '8.66'
$ws.Range("E45").Value = '  -1.95%  '

$ws.Range("E46").Value = '  -3.16%  '

Set-TextValue "D47" '406.63'
$ws.Range("E47").Value = '  -3.21%  '

Set-TextValue "D48" '45.53'
$ws.Range("E48").Value = '  -2.45%  '

Set-TextValue "D49" '143.51'
$ws.Range("E49").Value = '  +0.52%  '

Set-TextValue "D50" '0.000271'
$ws.Range("E50").Value = '  -10.84%  '

Set-TextValue "D51" '39.99'
$ws.Range("E51").Value = '  +4.02%  '
